# Apply "Add 2022-06-23 data" update to Fonds de solidarite volet 1 dataset.
# For a set of rows, update column C (nombre_aides) and column E (montant_total)
# to reflect newly aggregated values, leaving all other cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10;  C = 278207; E = 1752504854 },
    @{ Row = 111; C = 6012;   E = 12096389 },
    @{ Row = 115; C = 17559;  E = 38636494 },
    @{ Row = 117; C = 19729;  E = 56645610 },
    @{ Row = 124; C = 2663;   E = 3903857 },
    @{ Row = 125; C = 4606;   E = 13159157 },
    @{ Row = 126; C = 5646;   E = 8180829 },
    @{ Row = 134; C = 5683;   E = 17181109 },
    @{ Row = 139; C = 3321;   E = 9344392 },
    @{ Row = 168; C = 285089; E = 1212846471 },
    @{ Row = 169; C = 562653; E = 1285790763 },
    @{ Row = 170; C = 367525; E = 2847630590 },
    @{ Row = 171; C = 115211; E = 448574472 },
    @{ Row = 174; C = 357339; E = 1019792833 },
    @{ Row = 177; C = 96775;  E = 174802756 },
    @{ Row = 179; C = 235780; E = 813277362 },
    @{ Row = 180; C = 141516; E = 341192575 },
    @{ Row = 254; C = 62455;  E = 99573670 }
)

foreach ($u in $updates) {
    $ws.Range("C$($u.Row)").Value = $u.C
    $ws.Range("E$($u.Row)").Value = $u.E
}
